$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 317, shifting existing rows 317:430 down to 318:431
$ws.Rows.Item(317).Insert()

# Populate the newly inserted row 317 with the new record's data
$ws.Cells.Item(317, 1).Value = 5
$ws.Cells.Item(317, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(317, 3).Value = "Maule"
$ws.Cells.Item(317, 4).Value = 44900
$ws.Cells.Item(317, 5).Value = 7
$ws.Cells.Item(317, 6).Value = 100114013
$ws.Cells.Item(317, 7).Value = "Zanahoria"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 400
$ws.Cells.Item(317, 11).Value = 8000
$ws.Cells.Item(317, 12).Value = 8000
$ws.Cells.Item(317, 13).Value = 8000
$ws.Cells.Item(317, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(317, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(317, 16).Value = 400
$ws.Cells.Item(317, 17).Value = 20
$ws.Cells.Item(317, 18).Value = "Hortaliza"
